$d = $word.ActiveDocument

# The document has three "Poate edita/șterge ..." bullet points describing
# edit/delete operations ("un film", "lista de categorii", "un articol").
# Previously only the embedded word "șterge" (delete) was struck through in
# each; now the whole line in each of these three paragraphs must be struck
# through in full (the functionality is being marked as fully resolved /
# crossed out), and the runs that become uniformly-formatted collapse into a
# single run, same as Word does when it rewrites a paragraph's runs.
#
# Paragraph "...un film" keeps its trailing "un film" run separate from the
# rest (that run carries its own revision-save id in the original file), so
# only the leading "Poate edita/" + "șterge" + " " runs are merged together.
# Paragraphs "...lista de categorii" and "...un articol" collapse entirely
# into a single run each.

function Find-ParagraphByText($doc, $exactText) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -eq $exactText) {
            return $p
        }
    }
    return $null
}

function Strike-And-Merge($doc, $fullText, $leadLen, $leadText) {
    # 1) Strike the leading part (runs that must merge) - setting Font
    #    directly (not via Find/Replace) preserves existing run boundaries.
    $p = Find-ParagraphByText $doc $fullText
    $start = $p.Range.Start
    $lead = $doc.Range($start, $start + $leadLen)
    $lead.Font.StrikeThrough = 1

    # 2) Re-find the (now fully-struck) leading span and do a no-op
    #    Find/Replace over it: Word collapses a replaced span that ends up
    #    with uniform formatting into a single run.
    $p2 = Find-ParagraphByText $doc $fullText
    $start2 = $p2.Range.Start
    $lead2 = $doc.Range($start2, $start2 + $leadLen)
    $lead2.Find.Execute($leadText, $true, $false, $false, $false, $false, $true, 1, $false, $leadText, 2) | Out-Null
}

# --- Paragraph: "Poate edita/șterge un film" ---------------------------
Strike-And-Merge $d "Poate edita/șterge un film`r" 19 "Poate edita/șterge "
# Strike the trailing "un film" run too (kept as its own run).
$pFilm = Find-ParagraphByText $d "Poate edita/șterge un film`r"
$pFilm.Range.Font.StrikeThrough = 1

# --- Paragraph: "Poate edita/șterge lista de categorii" -----------------
$catText = "Poate edita/șterge lista de categorii"
Strike-And-Merge $d "$catText`r" $catText.Length $catText

# --- Paragraph: "Poate edita/șterge un articol" --------------------------
$artText = "Poate edita/șterge un articol"
Strike-And-Merge $d "$artText`r" $artText.Length $artText

Write-Output "done"
